$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Percentage cells (column H) need explicit Text format to avoid Excel
# auto-converting "NN%" strings into percentage numbers.
$ws.Range('H3').NumberFormat = "@"
$ws.Range('H8').NumberFormat = "@"
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H19').NumberFormat = "@"
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H27').NumberFormat = "@"
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H33').NumberFormat = "@"
$ws.Range('H39').NumberFormat = "@"
$ws.Range('H44').NumberFormat = "@"
$ws.Range('H46').NumberFormat = "@"

$ws.Range('E2').Value = '2026-02-27 04:18:35'
$ws.Range('N2').Value = '0.6 °C 3:57 TU'
$ws.Range('E3').Value = '2026-02-27 04:18:37'
$ws.Range('H3').Value = '34%'
$ws.Range('K3').Value = '-0.1 MJ/m2'
$ws.Range('E4').Value = '2026-02-27 04:18:39'
$ws.Range('J4').Value = '1026.1 hPa'
$ws.Range('N4').Value = '6.4 °C 3:58 TU'
$ws.Range('E5').Value = '2026-02-27 04:18:42'
$ws.Range('E6').Value = '2026-02-27 04:18:44'
$ws.Range('J6').Value = '1025.9 hPa'
$ws.Range('N6').Value = '9.0 °C 3:50 TU'
$ws.Range('E7').Value = '2026-02-27 04:18:47'
$ws.Range('J7').Value = '1026.2 hPa'
$ws.Range('N7').Value = '9.4 °C 3:59 TU'
$ws.Range('O7').Value = '10.4 °C'
$ws.Range('E8').Value = '2026-02-27 04:18:49'
$ws.Range('H8').Value = '44%'
$ws.Range('J8').Value = '1025.7 hPa'
$ws.Range('E9').Value = '2026-02-27 04:18:52'
$ws.Range('O9').Value = '8.2 °C'
$ws.Range('E10').Value = '2026-02-27 04:18:54'
$ws.Range('O10').Value = '9.3 °C'
$ws.Range('E11').Value = '2026-02-27 04:18:57'
$ws.Range('N11').Value = '1.7 °C 3:59 TU'
$ws.Range('O11').Value = '2.7 °C'
$ws.Range('E12').Value = '2026-02-27 04:18:59'
$ws.Range('M12').Value = '9.0 °C 3:58 TU'
$ws.Range('O12').Value = '7.6 °C'
$ws.Range('E13').Value = '2026-02-27 04:19:01'
$ws.Range('L13').Value = '5.4 km/h - 151º 3:37 TU'
$ws.Range('O13').Value = '-1.0 °C'
$ws.Range('E14').Value = '2026-02-27 04:19:04'
$ws.Range('L14').Value = '17.3 km/h - 325º 3:40 TU'
$ws.Range('N14').Value = '8.3 °C 3:59 TU'
$ws.Range('O14').Value = '9.4 °C'
$ws.Range('E15').Value = '2026-02-27 04:19:06'
$ws.Range('H15').Value = '96%'
$ws.Range('O15').Value = '8.1 °C'
$ws.Range('E16').Value = '2026-02-27 04:19:09'
$ws.Range('L16').Value = '20.2 km/h - 227º 3:54 TU'
$ws.Range('E17').Value = '2026-02-27 04:19:11'
$ws.Range('L17').Value = '40.3 km/h - 224º 3:32 TU'
$ws.Range('N17').Value = '6.9 °C 3:45 TU'
$ws.Range('E18').Value = '2026-02-27 04:19:14'
$ws.Range('J18').Value = '1026.0 hPa'
$ws.Range('N18').Value = '9.3 °C 3:56 TU'
$ws.Range('O18').Value = '10.0 °C'
$ws.Range('E19').Value = '2026-02-27 04:19:16'
$ws.Range('H19').Value = '85%'
$ws.Range('N19').Value = '7.5 °C 3:55 TU'
$ws.Range('O19').Value = '7.9 °C'
$ws.Range('E20').Value = '2026-02-27 04:19:19'
$ws.Range('O20').Value = '1.6 °C'
$ws.Range('E21').Value = '2026-02-27 04:19:21'
$ws.Range('H21').Value = '78%'
$ws.Range('N21').Value = '2.6 °C 3:53 TU'
$ws.Range('O21').Value = '4.2 °C'
$ws.Range('E22').Value = '2026-02-27 04:19:24'
$ws.Range('O22').Value = '0.8 °C'
$ws.Range('E23').Value = '2026-02-27 04:19:26'
$ws.Range('N23').Value = '2.1 °C 3:38 TU'
$ws.Range('O23').Value = '2.9 °C'
$ws.Range('E24').Value = '2026-02-27 04:19:29'
$ws.Range('J24').Value = '1026.4 hPa'
$ws.Range('N24').Value = '2.2 °C 3:49 TU'
$ws.Range('O24').Value = '5.2 °C'
$ws.Range('E25').Value = '2026-02-27 04:19:31'
$ws.Range('E26').Value = '2026-02-27 04:19:34'
$ws.Range('H26').Value = '47%'
$ws.Range('J26').Value = '1025.0 hPa'
$ws.Range('M26').Value = '8.1 °C 3:46 TU'
$ws.Range('O26').Value = '7.5 °C'
$ws.Range('E27').Value = '2026-02-27 04:19:37'
$ws.Range('H27').Value = '50%'
$ws.Range('L27').Value = '27.4 km/h - 238º 3:32 TU'
$ws.Range('E28').Value = '2026-02-27 04:19:39'
$ws.Range('J28').Value = '1026.2 hPa'
$ws.Range('N28').Value = '4.9 °C 3:58 TU'
$ws.Range('O28').Value = '5.9 °C'
$ws.Range('E29').Value = '2026-02-27 04:19:42'
$ws.Range('E30').Value = '2026-02-27 04:19:44'
$ws.Range('J30').Value = '1025.8 hPa'
$ws.Range('E31').Value = '2026-02-27 04:19:47'
$ws.Range('J31').Value = '1025.4 hPa'
$ws.Range('E32').Value = '2026-02-27 04:19:50'
$ws.Range('H32').Value = '93%'
$ws.Range('O32').Value = '1.5 °C'
$ws.Range('E33').Value = '2026-02-27 04:19:52'
$ws.Range('H33').Value = '69%'
$ws.Range('J33').Value = '1029.5 hPa'
$ws.Range('N33').Value = '1.5 °C 3:56 TU'
$ws.Range('O33').Value = '2.9 °C'
$ws.Range('E34').Value = '2026-02-27 04:19:55'
$ws.Range('L34').Value = '16.9 km/h - 21º 3:58 TU'
$ws.Range('O34').Value = '1.7 °C'
$ws.Range('E35').Value = '2026-02-27 04:19:57'
$ws.Range('J35').Value = '1025.6 hPa'
$ws.Range('E36').Value = '2026-02-27 04:20:00'
$ws.Range('J36').Value = '1026.1 hPa'
$ws.Range('M36').Value = '10.1 °C 3:43 TU'
$ws.Range('O36').Value = '8.9 °C'
$ws.Range('E37').Value = '2026-02-27 04:20:03'
$ws.Range('J37').Value = '1028.8 hPa'
$ws.Range('N37').Value = '2.0 °C 3:51 TU'
$ws.Range('O37').Value = '2.9 °C'
$ws.Range('E38').Value = '2026-02-27 04:20:05'
$ws.Range('L38').Value = '9.7 km/h - 282º 3:55 TU'
$ws.Range('N38').Value = '7.1 °C 3:58 TU'
$ws.Range('O38').Value = '7.9 °C'
$ws.Range('E39').Value = '2026-02-27 04:20:08'
$ws.Range('H39').Value = '19%'
$ws.Range('K39').Value = '-0.1 MJ/m2'
$ws.Range('L39').Value = '23.8 km/h - 303º 3:32 TU'
$ws.Range('M39').Value = '6.1 °C 3:43 TU'
$ws.Range('O39').Value = '5.2 °C'
$ws.Range('E40').Value = '2026-02-27 04:20:10'
$ws.Range('N40').Value = '1.3 °C 3:54 TU'
$ws.Range('O40').Value = '2.3 °C'
$ws.Range('E41').Value = '2026-02-27 04:20:13'
$ws.Range('J41').Value = '1026.2 hPa'
$ws.Range('N41').Value = '7.7 °C 3:59 TU'
$ws.Range('O41').Value = '9.3 °C'
$ws.Range('E42').Value = '2026-02-27 04:20:15'
$ws.Range('M42').Value = '9.2 °C 3:45 TU'
$ws.Range('O42').Value = '8.0 °C'
$ws.Range('E43').Value = '2026-02-27 04:20:18'
$ws.Range('N43').Value = '2.9 °C 3:59 TU'
$ws.Range('O43').Value = '4.4 °C'
$ws.Range('E44').Value = '2026-02-27 04:20:20'
$ws.Range('H44').Value = '69%'
$ws.Range('E45').Value = '2026-02-27 04:20:23'
$ws.Range('N45').Value = '5.1 °C 3:59 TU'
$ws.Range('O45').Value = '7.0 °C'
$ws.Range('E46').Value = '2026-02-27 04:20:26'
$ws.Range('H46').Value = '99%'
$ws.Range('J46').Value = '1026.4 hPa'
$ws.Range('N46').Value = '4.8 °C 3:57 TU'
$ws.Range('O46').Value = '7.4 °C'
